$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.031.95'
$ws.Range('E2').Value = '  -0.44%  '

$ws.Range('D3').Value = '1.829.28'
$ws.Range('E3').Value = '  -0.11%  '

$ws.Range('D4').Value = '''0.9989'

$ws.Range('D5').Value = '''241.05'
$ws.Range('E5').Value = '  -0.23%  '

$ws.Range('E6').Value = '  -5.49%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = '''0.07538'
$ws.Range('E8').Value = '  +1.83%  '

$ws.Range('D9').Value = '''44.59'
$ws.Range('E9').Value = '  +6.68%  '

$ws.Range('D10').Value = '''0.2908'
$ws.Range('E10').Value = '  -0.67%  '

$ws.Range('D11').Value = '''22.76'
$ws.Range('E11').Value = '  -0.65%  '

$ws.Range('E12').Value = '  -1.51%  '

$ws.Range('D13').Value = '1.831.43'
$ws.Range('E13').Value = '  +0.51%  '

$ws.Range('D14').Value = '''4.953'

$ws.Range('D15').Value = '''0.6639'
$ws.Range('E15').Value = '  -0.26%  '

$ws.Range('D16').Value = '''82.23'
$ws.Range('E16').Value = '  -0.63%  '

$ws.Range('D17').Value = '''0.000009105'
$ws.Range('E17').Value = '  +8.08%  '

$ws.Range('D18').Value = '''5.997'
$ws.Range('E18').Value = '  -1.95%  '

$ws.Range('D19').Value = '28.929.91'
$ws.Range('E19').Value = '  -0.74%  '

$ws.Range('D20').Value = '''224.39'

$ws.Range('D21').Value = '''12.32'
$ws.Range('E21').Value = '  -0.94%  '

$ws.Range('E22').Value = '  +0.02%  '

$ws.Range('D23').Value = '''7.182'
$ws.Range('E23').Value = '  +0.82%  '

$ws.Range('D24').Value = '''1.000'
$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('D25').Value = '''159.37'
$ws.Range('E25').Value = '  +0.37%  '

$ws.Range('D26').Value = '''8.381'
$ws.Range('E26').Value = '  -2.58%  '

$ws.Range('D27').Value = '''0.1354'
$ws.Range('E27').Value = '  -2.55%  '

$ws.Range('D28').Value = '''17.81'
$ws.Range('E28').Value = '  -0.50%  '

$ws.Range('E29').Value = '  -1.72%  '

$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '''4.048'
$ws.Range('E30').Value = '  -1.62%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '''4.026'
$ws.Range('E31').Value = '  -0.40%  '

$ws.Range('D32').Value = '''1.199'
$ws.Range('E32').Value = '  +1.00%  '

$ws.Range('D33').Value = '''0.05189'
$ws.Range('E33').Value = '  -1.23%  '

$ws.Range('D34').Value = '''1.834'
$ws.Range('E34').Value = '  -1.47%  '

$ws.Range('D35').Value = '''1.151'
$ws.Range('E35').Value = '  +0.86%  '

$ws.Range('D36').Value = '''0.7309'
$ws.Range('E36').Value = '  -1.30%  '

$ws.Range('D37').Value = '''2.610'
$ws.Range('E37').Value = '  -1.59%  '

$ws.Range('D38').Value = '1.286.72'
$ws.Range('E38').Value = '  -1.22%  '

$ws.Range('D39').Value = '''2.756'
$ws.Range('E39').Value = '  +0.91%  '

$ws.Range('D40').Value = '''0.01777'
$ws.Range('E40').Value = '  -0.69%  '

$ws.Range('D41').Value = '''6.392'
$ws.Range('E41').Value = '  +7.43%  '

$ws.Range('D42').Value = '''0.8926'
$ws.Range('E42').Value = '  -3.19%  '

$ws.Range('E43').Value = '  +0.16%  '

$ws.Range('D44').Value = '''101.51'
$ws.Range('E44').Value = '  -0.63%  '

$ws.Range('D45').Value = '1.980.35'
$ws.Range('E45').Value = '  +0.21%  '

$ws.Range('D47').Value = '''63.48'
$ws.Range('E47').Value = '  +0.30%  '

$ws.Range('E48').Value = '  -0.47%  '

$ws.Range('D49').Value = '''0.3969'
$ws.Range('E49').Value = '  -0.87%  '

$ws.Range('D50').Value = '''8.868'
$ws.Range('E50').Value = '  +1.35%  '

$ws.Range('D51').Value = '''1.647'
$ws.Range('E51').Value = '  -5.94%  '
